# Adds bite-rate-per-wave-energy data for the Bejarano paper (James), with
# low/medium/high wave conditions, to the herb bite rate estimates sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: row, wave-condition label (col B), bite-rate value (col D), species (col E)
$rows = @(
    @(2,  "low wave",    6.25,               "Calotomus carolinus"),
    @(3,  "low wave",    8.6111111111111107, "Siganus argenteus"),
    @(4,  "low wave",    19.2053264604811,   "Zebrasoma veliferum"),
    @(5,  "low wave",    2.9660652920962201, "Acanthurus pyroferus"),
    @(6,  "low wave",    10.95505617977528,  "Chlorurus microrhinos"),
    @(7,  "low wave",    37.604215048377434, "Scarus dimidiatus"),
    @(8,  "low wave",    1.1235955056179774, "Scarus oviceps"),
    @(9,  "low wave",    2.7083333333333335, "Siganus corallinus"),
    @(10, "medium wave", 35.799934832192896, "Siganus corallinus"),
    @(11, "low wave",    31.365629984051036, "Naso lituratus"),
    @(12, "medium wave", 33.463622776152754, "Naso lituratus"),
    @(13, "high wave",   7.2330934173039436, "Naso lituratus"),
    @(14, "low wave",    328.12092005062789, "Zebrasoma scopas"),
    @(15, "medium wave", 253.73992677916567, "Zebrasoma scopas"),
    @(16, "high wave",   6.4016064257028118, "Zebrasoma scopas"),
    @(17, "low wave",    18.166035353535349, "Scarus chameleon"),
    @(18, "medium wave", 21.521416083916083, "Scarus chameleon"),
    @(19, "high wave",   10.464456391875746, "Scarus chameleon"),
    @(20, "low wave",    0.55555555555555558,"Scarus forsteni"),
    @(21, "medium wave", 4.7422680412371134, "Scarus forsteni"),
    @(22, "high wave",   9.3862007168458792, "Scarus forsteni"),
    @(23, "low wave",    40.889554611745936, "Scarus niger"),
    @(24, "medium wave", 29.112461647740066, "Scarus niger"),
    @(25, "high wave",   18.765962307496867, "Scarus niger")
)

# The workbook's shared-string table introduces "medium wave", "high wave"
# and "low wave" (in that order) before any other new content, matching how
# the data was originally entered. Seed that order explicitly first.
$ws.Cells.Item(10, 2).Value = "medium wave"
$ws.Cells.Item(13, 2).Value = "high wave"
$ws.Cells.Item(2, 2).Value = "low wave"

foreach ($row in $rows) {
    $r = $row[0]
    $wave = $row[1]
    $value = $row[2]
    $species = $row[3]

    $ws.Cells.Item($r, 1).Value = "Bejarano"
    $ws.Cells.Item($r, 2).Value = $wave
    $ws.Cells.Item($r, 4).Value = $value
    $ws.Cells.Item($r, 5).Value = $species
    $ws.Cells.Item($r, 6).Value = "James"
}

# Center the bite-rate column, matching a couple of stray selections that
# picked up neighboring cells along the way.
$ws.Range("D2:D25").HorizontalAlignment = -4108
$ws.Range("C9:C10").HorizontalAlignment = -4108
$ws.Range("I21:I23").HorizontalAlignment = -4108

$ws.Range("D19").Select()
